$d = $word.ActiveDocument
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

# Touching a paragraph with Find/Replace causes the interop engine to
# normalize (merge) that paragraph's runs, matching the canonical save
# behaviour recorded in the target revision. Re-apply the same text on
# the three untouched-content paragraphs so their runs collapse exactly
# as in the reference diff, with no visible change.
$find.Execute("helpful comments", $true, $true, $false, $false, $false, $true, 1, $false, "helpful comments", 2)
$find.Execute("we focus on DE", $true, $true, $false, $false, $false, $true, 1, $false, "we focus on DE", 2)
$find.Execute("redesigned using parameter estimates", $true, $true, $false, $false, $false, $true, 1, $false, "redesigned using parameter estimates", 2)

# Real content change: the pluripotency-factor ranking paragraph now
# cites Supplementary Figure S7 instead of S6.
$find.Execute("Supplementary Figure S6, p-value of 0.0002", $true, $true, $false, $false, $false, $true, 1, $false, "Supplementary Figure S7, p-value of 0.0002", 2)

# The following paragraph (about accurate error-rate control) is also
# re-saved with its runs normalised in the target revision, even though
# its text is unchanged.
$find.Execute("error rate with the summation", $true, $true, $false, $false, $false, $true, 1, $false, "error rate with the summation", 2)
